$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.621.53"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "2.639.61"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("E4").Value = "  -0.07%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "607.65"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.38%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "148.17"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.29%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("E9").Value = "  +2.82%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "5.59"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.58%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.375"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.86%  "

$ws.Range("E12").Value = "  +0.01%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "27.66"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").Value = "3.109.43"
$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").Value = "63.463.29"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("E16").Value = "  +2.90%  "

$ws.Range("D17").Value = "2.629.46"
$ws.Range("E17").Value = "  -0.26%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "11.61"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.02%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.59"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +5.14%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "345.64"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.99%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.90"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.85%  "

$ws.Range("E22").Value = "  -0.21%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.56"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.69%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "67.04"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.21%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.72"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.04%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.12"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +8.24%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.59"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.13%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "554.21"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +6.17%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.163"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.12%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.02"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.54%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "

$ws.Range("E32").Value = "  +4.31%  "

$ws.Range("D33").Value = "0.0₃0859"
$ws.Range("E33").Value = "  +6.45%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.77"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.00%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.21"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +5.73%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "167.36"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -4.07%  "

$ws.Range("E37").Value = "  +1.37%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.25%  "

$ws.Range("E39").Value = "  +8.81%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "19.18"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.86%  "

$ws.Range("E41").Value = "  +0.04%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "166.14"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.17%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.81"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.27%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "22.20"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.35%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0574"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.81%  "

$ws.Range("E46").Value = "  +0.19%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0248"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.36%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0965"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.95"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +14.20%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "18.92"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.46%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.185"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +7.05%  "
